$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.830.08"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.95"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.42"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5030"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06409"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.72"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07697"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.253"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.31"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.859.39"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7923"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.60"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.850.93"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.37"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.327"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.928"
$ws.Range("E25").Value = "  +9.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.29"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.710"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05022"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.241"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.188"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.538"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.354"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174.23"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8921"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5590"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01562"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.550"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.667"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8080"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.46"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.771.49"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4513"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.92"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05038"
$ws.Range("E51").Value = "  -0.45%  "
